$d = $word.ActiveDocument

# 1. Sequencer update: "NovaSeq 6000" -> "NovaSeq X Plus (Australian Genome Research Facility)"
$d.Content.Find.Execute(
    "NovaSeq 6000",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "NovaSeq X Plus (Australian Genome Research Facility)",
    2
)

# 2. Report date field text update: "16-Sep-2024" -> "4-Mar-2025"
$d.Content.Find.Execute(
    "16-Sep-2024",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "4-Mar-2025",
    2
)
